# Daily attendance processing - 2025-11-24 08:33:56
# Applies the "Recorded By" list re-ordering / additions, the derived
# stat updates that follow from them, and the session 27 status flip
# (Pending -> Not Recorded) to the Y2 B2526 GIT & Liver attendance sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Recorded By" (column G) list updates -------------------------------

$ws.Range("G2").Value = "System, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G3").Value = "majorelle.magdy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, System, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G4").Value = "majorelle.magdy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg"
$ws.Range("G5").Value = "asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G7").Value = "Kerelos.zareef@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Fatmaelhady@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg"
$ws.Range("G9").Value = "Safa.hany@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg"
$ws.Range("G12").Value = "yassmina.fattoh@med.asu.edu.eg, dina.adel@med.asu.edu.eg"
$ws.Range("G15").Value = "Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("G28").Value = "maryam.ashraf@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg"

# --- Students recorded (column H) -----------------------------------------

$ws.Range("H12").Value = "37/251"

# --- Session 7 (row 7/8 of the BIOCHEMISTRY LAB/CBL group) is now recorded,
#     shifting it from "Missing" to a counted session, so the summary
#     Missing/Pending counters move by one in both the top Class Statistics
#     block (K/L) and the Group Statistics row (row 15, columns P/Q).

$ws.Range("L7").Value = 2
$ws.Range("L8").Value = 15
$ws.Range("P15").Value = 2
$ws.Range("Q15").Value = 15

# --- Average attendance % moved from 23.6% to 24.1% (Class Statistics L10
#     and Group Statistics S15). Force literal text so Excel does not
#     reinterpret "24.1%" as a numeric percentage.

$ws.Range("L10").Value = "'24.1%"
$ws.Range("S15").Value = "'24.1%"

# --- Row 27 (PHARMACOLOGY session 2) flips from "Pending" to
#     "Not Recorded" -- copy the formatting already used for that status
#     (row 29) onto row 27, then update the status text.

$src = $ws.Range("A29:I29")
$dst = $ws.Range("A27:I27")
$src.Copy()
$dst.PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("I27").Value = "Not Recorded"
